$wb = $excel.ActiveWorkbook

# --- Suite3 sheet: mark 8 more scripts as created (N -> Y) ---
$suite3 = $wb.Worksheets.Item("Suite3")
$suite3.Activate() | Out-Null

$rowsToFlip = @(6,7,8,13,14,15,16,17)
foreach ($r in $rowsToFlip) {
    $suite3.Cells.Item($r, 2).Value = "Y"
}

# --- AppControl sheet: selection moves to E5 (the Suite3 pass-count cell) ---
# Selecting a range on another sheet activates that sheet, so do this before
# re-activating Suite3 below (Suite3 is the tab that should end up selected).
$appControl = $wb.Worksheets.Item("AppControl")
$appControl.Range("E5").Select() | Out-Null

# --- Back to Suite3, which remains the active/selected tab ---
# Selection ends up on B18 with the view scrolled down so row 4 is at the top
$suite3.Activate() | Out-Null
$suite3.Range("B18").Select() | Out-Null
